$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix up species / row labels: strip the long trailing underscore "fill-in-the-blank"
# runs from the scanned-form labels, fix typos, and fill in the label that was
# left blank (row 5 -> "Cabezon"). Order matters for shared-string layout parity.
$ws.Range("A4").Value = "Bonito Pacific"
$ws.Range("A6").Value = "Halibut California"
$ws.Range("A8").Value = "Mackerel jack"
$ws.Range("A14").Value = "Sheepbead California"
$ws.Range("A5").Value = "Cabezon"
$ws.Range("A20").Value = "Total number of fish"
$ws.Range("A22").Value = "Number of angler days"
$ws.Range("A12").Value = "Sculpin"
$ws.Range("A15").Value = "Tuna albacore"
$ws.Range("A18").Value = "Yellowtail California"

# Those rows no longer need the extra wrapped height the underscores forced.
$ws.Rows.Item(4).RowHeight = 17
$ws.Rows.Item(6).RowHeight = 17
$ws.Rows.Item(8).RowHeight = 17
$ws.Rows.Item(14).RowHeight = 17
$ws.Rows.Item(15).RowHeight = 17
$ws.Rows.Item(18).RowHeight = 17
$ws.Rows.Item(20).RowHeight = 17
$ws.Rows.Item(22).RowHeight = 17

# Scroll back to the top and leave the cursor on A19.
[void]$ws.Range("A19").Select()
